# Apply cryptos list update (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.643.05'
$ws.Range('E2').Value = '  -1.40%  '

$ws.Range('D3').Value = '2.778.19'
$ws.Range('E3').Value = '  -2.25%  '

$ws.Range('D5').Value = "'" + '359.22'
$ws.Range('E5').Value = '  -0.63%  '

$ws.Range('D6').Value = "'" + '109.02'
$ws.Range('E6').Value = '  -3.41%  '

$ws.Range('D7').Value = "'" + '0.556'
$ws.Range('E7').Value = '  -2.85%  '

$ws.Range('E8').Value = '  +0.18%  '

$ws.Range('D9').Value = "'" + '0.585'
$ws.Range('E9').Value = '  -3.94%  '

$ws.Range('D10').Value = "'" + '39.66'
$ws.Range('E10').Value = '  -3.78%  '

$ws.Range('E11').Value = '  +4.34%  '

$ws.Range('D12').Value = "'" + '0.0843'
$ws.Range('E12').Value = '  -3.21%  '

$ws.Range('D13').Value = "'" + '19.68'
$ws.Range('E13').Value = '  -2.36%  '

$ws.Range('E14').Value = '  -2.90%  '

$ws.Range('D15').Value = '3.219.72'
$ws.Range('E15').Value = '  -2.08%  '

$ws.Range('D16').Value = '2.776.63'
$ws.Range('E16').Value = '  -3.14%  '

$ws.Range('D17').Value = "'" + '0.926'
$ws.Range('E17').Value = '  -1.03%  '

$ws.Range('D18').Value = '51.641.41'
$ws.Range('E18').Value = '  -1.24%  '

$ws.Range('D19').Value = "'" + '7.66'
$ws.Range('E19').Value = '  +1.06%  '

$ws.Range('E20').Value = '  -2.03%  '

$ws.Range('D21').Value = "'" + '13.19'
$ws.Range('E21').Value = '  -2.44%  '

$ws.Range('D22').Value = '0.0₃0968'
$ws.Range('E22').Value = '  -3.31%  '

$ws.Range('D23').Value = "'" + '70.00'
$ws.Range('E23').Value = '  -0.92%  '

$ws.Range('D24').Value = "'" + '267.98'
$ws.Range('E24').Value = '  -1.87%  '

$ws.Range('D25').Value = "'" + '2.77'
$ws.Range('E25').Value = '  -2.58%  '

$ws.Range('D26').Value = "'" + '26.32'
$ws.Range('E26').Value = '  -2.61%  '

$ws.Range('D27').Value = "'" + '0.999'
$ws.Range('E27').Value = '  -0.04%  '

$ws.Range('D28').Value = "'" + '0.165'
$ws.Range('E28').Value = '  +14.61%  '

$ws.Range('D29').Value = "'" + '10.16'
$ws.Range('E29').Value = '  -2.03%  '

$ws.Range('D30').Value = "'" + '2.26'
$ws.Range('E30').Value = '  +0.23%  '

$ws.Range('D31').Value = "'" + '35.30'
$ws.Range('E31').Value = '  -0.67%  '

$ws.Range('D32').Value = "'" + '6.14'
$ws.Range('E32').Value = '  +3.58%  '

$ws.Range('D33').Value = "'" + '52.05'
$ws.Range('E33').Value = '  -0.87%  '

$ws.Range('D34').Value = "'" + '0.0443'
$ws.Range('E34').Value = '  -8.86%  '

$ws.Range('D35').Value = "'" + '0.0839'
$ws.Range('E35').Value = '  -1.97%  '

$ws.Range('D36').Value = "'" + '5.16'

$ws.Range('D37').Value = "'" + '1.00'
$ws.Range('E37').Value = '  -0.03%  '

$ws.Range('D38').Value = "'" + '18.80'
$ws.Range('E38').Value = '  +1.31%  '

$ws.Range('E39').Value = '  -5.08%  '

$ws.Range('D40').Value = "'" + '1.95'
$ws.Range('E40').Value = '  -5.00%  '

$ws.Range('E41').Value = '  -3.35%  '

$ws.Range('D42').Value = "'" + '2.50'
$ws.Range('E42').Value = '  -1.69%  '

$ws.Range('E43').Value = '  -3.71%  '

$ws.Range('D44').Value = "'" + '119.74'
$ws.Range('E44').Value = '  -5.84%  '

$ws.Range('D45').Value = "'" + '21.72'
$ws.Range('E45').Value = '  -5.83%  '

$ws.Range('D46').Value = '2.090.54'
$ws.Range('E46').Value = '  -0.36%  '

$ws.Range('D47').Value = "'" + '3.25'
$ws.Range('E47').Value = '  -3.34%  '

$ws.Range('D48').Value = "'" + '2.30'
$ws.Range('E48').Value = '  +0.34%  '

$ws.Range('D49').Value = "'" + '0.932'
$ws.Range('E49').Value = '  -4.74%  '

$ws.Range('D50').Value = "'" + '5.55'
$ws.Range('E50').Value = '  -6.45%  '

$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = "'" + '1.28'
$ws.Range('E51').Value = '  +3.45%  '
